$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first occurrence of the data rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 30
$ws1.Range("F5").Value = 4846
$ws1.Range("F7").Value = 72
$ws1.Range("F8").Value = 271
$ws1.Range("F9").Value = 37

# Sheet "全部类型" (All types) - second occurrence of the same data rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 30
$ws4.Range("F9").Value = 4846
$ws4.Range("F11").Value = 72
$ws4.Range("F13").Value = 271
$ws4.Range("F14").Value = 37
